$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3315
$ws.Range("C2").Value = 4

$ws.Range("A3").Value = 9821

$ws.Range("A4").Value = 6117
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = 6541
$ws.Range("C5").Value = 5

$ws.Range("A6").Value = 5724
$ws.Range("C6").Value = 4
